# Fixed bugs: Not being able to show status result of bids after round 1.
# Updated BM under 27

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Log")
$ws.Activate()

# --- Update existing "Iteration" values for rows 27 and 28 (3 -> 4) ---
$ws.Range("B27").Value = 4
$ws.Range("B28").Value = 4

# --- Add the new bug-log entry in row 29 ---
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = 4
$ws.Range("C29").Value = "Landing page "
$ws.Range("D29").Value = "Not being able to show status result of bids after round 1"
$ws.Range("E29").Value = "Resolved"
$ws.Range("F29").Value = "'14/11/2019"
$ws.Range("G29").Value = "'14/11/2019"
$ws.Range("H29").Value = "Matthew & DaEun"

# --- Match the formatting used by the preceding data row (27) ---
$ws.Range("A27:H27").Copy() | Out-Null
$ws.Range("A29:H29").PasteSpecial(-4122) | Out-Null

# Column A on row 29 uses the "interior" border variant (same as column C)
# rather than continuing the alternating outer-border pattern of A27/A28.
$ws.Range("C27").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null

# --- Reflect where the author left the viewport/selection ---
$ws.Range("D29").Select() | Out-Null
